$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new value. NumberFormat is forced to
# Text ("@") before assignment so numeric-looking strings (e.g. "1.000",
# "5.680", "0.00001095") keep their exact original text rather than being
# auto-coerced to a Double and losing trailing zeros / switching to
# scientific notation.
$updates = @{
    'D2' = '28.473.94'
    'E2' = '  +1.29%  '
    'D3' = '1.823.53'
    'E3' = '  +1.64%  '
    'D4' = '1.001'
    'E4' = '  +0.04%  '
    'D5' = '316.96'
    'E5' = '  +0.01%  '
    'E6' = '  +0.05%  '
    'D7' = '0.5412'
    'E7' = '  +0.96%  '
    'D8' = '0.4035'
    'E8' = '  +7.05%  '
    'D9' = '0.07669'
    'E9' = '  +2.65%  '
    'D10' = '1.121'
    'E10' = '  +2.46%  '
    'E11' = '  +0.49%  '
    'D12' = '6.329'
    'E12' = '  +3.64%  '
    'D13' = '7.642'
    'E13' = '  +6.00%  '
    'E14' = '  +0.02%  '
    'D15' = '20.93'
    'E15' = '  +1.47%  '
    'D16' = '1.823.77'
    'E16' = '  +2.57%  '
    'D17' = '0.00001095'
    'E17' = '  +3.56%  '
    'D18' = '89.84'
    'E18' = '  +0.77%  '
    'D19' = '0.06608'
    'E19' = '  +2.37%  '
    'D20' = '17.69'
    'E20' = '  +2.05%  '
    'E21' = '  +0.10%  '
    'D22' = '6.065'
    'E22' = '  +2.73%  '
    'D23' = '28.473.79'
    'E23' = '  +1.19%  '
    'D24' = '11.13'
    'D25' = '2.269'
    'E25' = '  +8.24%  '
    'D26' = '2.466'
    'E26' = '  +8.50%  '
    'B27' = 'Monero'
    'C27' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D27' = '157.84'
    'E27' = '  +1.96%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '20.74'
    'E28' = '  +2.55%  '
    'D29' = '2.034.69'
    'E29' = '  +2.46%  '
    'D30' = '123.87'
    'E30' = '  +2.67%  '
    'B31' = 'ImmutableX'
    'C31' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D31' = '1.122'
    'E31' = '  +0.43%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '0.1108'
    'E32' = '  +4.94%  '
    'D33' = '5.680'
    'E33' = '  +2.12%  '
    'B34' = 'HuobiToken'
    'C34' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D34' = '3.641'
    'E34' = '  -0.42%  '
    'B35' = 'Hedera'
    'C35' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D35' = '0.07346'
    'E35' = '  +12.11%  '
    'D36' = '0.2239'
    'E36' = '  -0.79%  '
    'D37' = '0.02344'
    'E37' = '  +2.65%  '
    'D38' = '5.209'
    'E38' = '  +3.66%  '
    'D39' = '8.848'
    'D40' = '11.36'
    'E40' = '  +2.49%  '
    'D41' = '0.6287'
    'E41' = '  +1.97%  '
    'D42' = '1.187'
    'E42' = '  +1.14%  '
    'D43' = '1.000'
    'E43' = '  +0.07%  '
    'D44' = '1.402'
    'E44' = '  -3.11%  '
    'D45' = '13.52'
    'E45' = '  +2.09%  '
    'D46' = '3.700'
    'E46' = '  +0.69%  '
    'D47' = '0.5851'
    'E47' = '  +1.23%  '
    'D48' = '125.27'
    'E48' = '  -1.68%  '
    'D49' = '2.003'
    'E49' = '  +3.93%  '
    'D50' = '1.198'
    'E50' = '  +0.71%  '
    'D51' = '0.06875'
    'E51' = '  +0.89%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
